{"js": "// Fix a units typo in the recitation worksheet: the density of air should\n// be given in g/cm^3 (grams per cubic centimeter), not g/dm^3.\n//\n// \"Calculate the mass of the air contained in a room ... given that the\n//  density of air is 1.29 g/dm^3 at 25 \u00b0C.\"\n//          -> \"... density of air is 1.29 g/cm^3 at 25 \u00b0C.\"\n\nconst body = context.document.body;\n\n// Locate the exact (unique) occurrence of the incorrect unit text.\nconst results = body.search(\"g/dm^3\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"g/cm^3\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Fix a units typo in the recitation worksheet: the density of air should\n# be given in g/cm^3 (grams per cubic centimeter), not g/dm^3.\n#\n# \"Calculate the mass of the air contained in a room ... given that the\n#  density of air is 1.29 g/dm^3 at 25 \u00b0C.\"\n#          -> \"... density of air is 1.29 g/cm^3 at 25 \u00b0C.\"\n\n$d = $word.ActiveDocument\n\n$needle = \"g/dm^3\"\n$replacement = \"g/cm^3\"\n\nforeach ($p in $d.Paragraphs) {\n    $rng = $p.Range\n    if ($rng.Text -like \"*$needle*\") {\n        # Paragraph.Range includes the trailing paragraph mark; trim it off\n        # before rewriting the text so we don't split the paragraph in two.\n        $rng.MoveEnd(1, -1) | Out-Null\n        $rng.Text = $rng.Text.Replace($needle, $replacement)\n        break\n    }\n}\n"}
